$wb = $excel.ActiveWorkbook

# StatOutput sheet: update the result row (row 2) with the new counts, forcing text so they remain
# shared-string cells like the rest of the sheet instead of being auto-converted to numbers.
$statWs = $wb.Worksheets.Item("StatOutput")
$statRow2 = $statWs.Range("A2:D2")
$statRow2.NumberFormat = "@"
$statWs.Range("A2").Value = "2"
$statWs.Range("B2").Value = "6"
$statWs.Range("C2").Value = "9"
$statWs.Range("D2").Value = "2"

# StatOutput_Message sheet: update the Cypher query text (row 18) to reference Boxer instead of Akita
$msgWs = $wb.Worksheets.Item("StatOutput_Message")
$boxerQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Boxer']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$msgWs.Range("A18").Value = $boxerQuery
